# Modified data generation parameters to get reasonable parameters.
#
# SimParameters sheet holds the relative-risk multipliers that drive the
# treated-pregnancy / treated-preeclampsia outcome sheets via formulas such
# as "=potential_preg_untrt!C9*SimParameters!$B$4". Swap the Low-severity and
# High-severity multipliers for both the Abortion block (B4/B6) and the
# Preeclampsia block (B8/B10); everything downstream recalculates on its own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SimParameters")

# Abortion severity multipliers
$ws.Range("B4").Value = 0.25   # Low severity      (was 0.75)
$ws.Range("B6").Value = 0.75   # High severity      (was 0.25)

# Preeclampsia severity multipliers
$ws.Range("B8").Value = 0.25   # Low severity      (was 0.75)
$ws.Range("B10").Value = 0.75  # High severity      (was 0.25)

# Make SimParameters the active sheet/tab, with cell B11 selected, matching
# the saved view state of the workbook after the edit.
$ws.Activate()
$ws.Range("B11").Select()
